$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# ExternalAct.classCode row (row 12): Min and Base Min go from 1 to 0
# (force text format first so Excel keeps "0" as text, matching the
# original shared-string cell type instead of converting it to a number)
$elem.Range("F12").NumberFormat = "@"
$elem.Range("F12").Value = "0"
$elem.Range("AG12").NumberFormat = "@"
$elem.Range("AG12").Value = "0"

# Binding Value Set URLs updated to CDA core value sets
$elem.Range("Z12").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActClass"
$elem.Range("Z13").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActMood"
